# pedidos.xlsx edit:
#  - delete the "megatronic" and "Pedidos" sheets entirely
#  - the old "PINI" sheet becomes "papa" (its row 2 product/quantity updated)
#  - the old "ooo" sheet (just headers) is renamed to "PINI"

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove sheets that disappear completely.
$wb.Worksheets.Item("megatronic").Delete()
$wb.Worksheets.Item("Pedidos").Delete()

# Free up the "PINI" name by temporarily renaming "ooo" (it will become
# the new "PINI" sheet further down, keeping its header-only content).
$wb.Worksheets.Item("ooo").Name = "PINI_TMP"

# The old "PINI" sheet becomes "papa"; update its second row of data.
$papa = $wb.Worksheets.Item("PINI")
$papa.Name = "papa"
$papa.Range("A2").Value = "Auricular Inalambrico"
# Force the quantity to be stored as text (matching the sheet's existing
# text-typed cells) instead of a native number, then drop back to the
# default style so no extra formatting sticks to the cell.
$papa.Range("B2").NumberFormat = "@"
$papa.Range("B2").Value = "5"
$papa.Range("B2").Style = "Normal"

# Finish the rename so the former "ooo" sheet is now called "PINI".
$wb.Worksheets.Item("PINI_TMP").Name = "PINI"
